# Regenerate merged AHB files:
#  - rename the "_old"/"_new" header-suffix columns to "_FV2310"/"_FV2404"
#  - freeze the header row
#  - turn the data range into a native Excel table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row shared strings ------------------------------
# Columns A..J carry the "..._old" suffix, columns L..U carry "..._new".
# (Column K holds the constant "diff" header and is left untouched.)
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$col = 0
foreach ($name in $baseNames) {
    $col = $col + 1
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2310"
}

$col = 11
foreach ($name in $baseNames) {
    $col = $col + 1
    $ws.Cells.Item(1, $col).Value = "$($name)_FV2404"
}

# --- 2. Freeze the header row ------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into a table -----------------------------
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U79"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Output "done"
